$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (shifts existing D:K to F:M)
$ws.Range("D:E").Insert()

# The newly inserted columns inherit formatting from the column to the left (C).
# Copy number/date formatting from column F (the old column D, now shifted) onto D:E
# so the new columns match the rest of the data block.
$ws.Range("F7:F102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$ws.Range("F7:F102").Copy()
$ws.Range("E7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Approximate the column widths of the two new columns off the neighbouring
# (already best-fit) column F, since the data has the same shape.
$ws.Columns("D").ColumnWidth = $ws.Columns("F").ColumnWidth()
$ws.Columns("E").ColumnWidth = $ws.Columns("F").ColumnWidth()

# Populate the new columns (D, E) with the latest two reporting periods, and
# refresh column F where the prior most-recent period was restated.

$ws.Cells.Item(7, 4).Value = 43465
$ws.Cells.Item(7, 5).Value = 43373
$ws.Cells.Item(7, 6).Value = 43281
$ws.Cells.Item(8, 4).Value = 692100
$ws.Cells.Item(8, 5).Value = 657600
$ws.Cells.Item(8, 6).Value = 617500
$ws.Cells.Item(9, 4).Value = 149700
$ws.Cells.Item(9, 5).Value = 133100
$ws.Cells.Item(9, 6).Value = 136600
$ws.Cells.Item(10, 4).Value = 542400
$ws.Cells.Item(10, 5).Value = 524500
$ws.Cells.Item(10, 6).Value = 480900
$ws.Cells.Item(12, 4).Value = 34100
$ws.Cells.Item(12, 5).Value = 21800
$ws.Cells.Item(12, 6).Value = 19100
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(13, 5).Value = 0
$ws.Cells.Item(13, 6).Value = 0
$ws.Cells.Item(14, 4).Value = 32000
$ws.Cells.Item(14, 5).Value = 26400
$ws.Cells.Item(14, 6).Value = 12900
$ws.Cells.Item(15, 4).Value = 261300
$ws.Cells.Item(15, 5).Value = 241800
$ws.Cells.Item(15, 6).Value = 238000
$ws.Cells.Item(17, 4).Value = 542200
$ws.Cells.Item(17, 5).Value = 446800
$ws.Cells.Item(17, 6).Value = 475400
$ws.Cells.Item(18, 4).Value = 149900
$ws.Cells.Item(18, 5).Value = 210800
$ws.Cells.Item(18, 6).Value = 142100
$ws.Cells.Item(20, 4).Value = 5900
$ws.Cells.Item(20, 5).Value = -20200
$ws.Cells.Item(20, 6).Value = -15800
$ws.Cells.Item(21, 4).Value = 417000
$ws.Cells.Item(21, 5).Value = 432400
$ws.Cells.Item(21, 6).Value = 364200
$ws.Cells.Item(22, 4).Value = 49500
$ws.Cells.Item(22, 5).Value = 43800
$ws.Cells.Item(22, 6).Value = 43900
$ws.Cells.Item(23, 4).Value = 106200
$ws.Cells.Item(23, 5).Value = 146800
$ws.Cells.Item(23, 6).Value = 82300
$ws.Cells.Item(24, 4).Value = 9200
$ws.Cells.Item(24, 5).Value = 51000
$ws.Cells.Item(24, 6).Value = 36400
$ws.Cells.Item(25, 4).Value = 0
$ws.Cells.Item(25, 5).Value = 0
$ws.Cells.Item(25, 6).Value = 0
$ws.Cells.Item(26, 4).Value = 96900
$ws.Cells.Item(26, 5).Value = 95800
$ws.Cells.Item(26, 6).Value = 45900
$ws.Cells.Item(27, 4).Value = 88600
$ws.Cells.Item(27, 5).Value = 95800
$ws.Cells.Item(27, 6).Value = 45900
$ws.Cells.Item(28, 4).Value = 0
$ws.Cells.Item(28, 5).Value = 0
$ws.Cells.Item(28, 6).Value = 0
$ws.Cells.Item(29, 4).Value = 14800
$ws.Cells.Item(29, 5).Value = -1800
$ws.Cells.Item(29, 6).Value = -400
$ws.Cells.Item(30, 4).Value = 0
$ws.Cells.Item(30, 5).Value = 0
$ws.Cells.Item(30, 6).Value = 0
$ws.Cells.Item(31, 4).Value = 0
$ws.Cells.Item(31, 5).Value = 0
$ws.Cells.Item(31, 6).Value = 0
$ws.Cells.Item(32, 4).Value = -5900
$ws.Cells.Item(32, 5).Value = 20200
$ws.Cells.Item(32, 6).Value = 15800
$ws.Cells.Item(33, 4).Value = 103400
$ws.Cells.Item(33, 5).Value = 93900
$ws.Cells.Item(33, 6).Value = 45500
$ws.Cells.Item(34, 4).Value = 0
$ws.Cells.Item(34, 5).Value = 0
$ws.Cells.Item(34, 6).Value = 0
$ws.Cells.Item(35, 4).Value = 103400
$ws.Cells.Item(35, 5).Value = 93900
$ws.Cells.Item(35, 6).Value = 45500
$ws.Cells.Item(38, 4).Value = 43465
$ws.Cells.Item(38, 5).Value = 43373
$ws.Cells.Item(38, 6).Value = 43281
$ws.Cells.Item(41, 4).Value = 387400
$ws.Cells.Item(41, 5).Value = 947700
$ws.Cells.Item(41, 6).Value = 901300
$ws.Cells.Item(42, 4).Value = "NA"
$ws.Cells.Item(42, 5).Value = "NA"
$ws.Cells.Item(42, 6).Value = "NA"
$ws.Cells.Item(43, 4).Value = 331900
$ws.Cells.Item(43, 5).Value = 274200
$ws.Cells.Item(43, 6).Value = 258400
$ws.Cells.Item(44, 4).Value = 87900
$ws.Cells.Item(44, 5).Value = 94600
$ws.Cells.Item(44, 6).Value = 91400
$ws.Cells.Item(45, 4).Value = 72700
$ws.Cells.Item(45, 5).Value = 64700
$ws.Cells.Item(45, 6).Value = 64300
$ws.Cells.Item(46, 4).Value = 879800
$ws.Cells.Item(46, 5).Value = 1381300
$ws.Cells.Item(46, 6).Value = 1315400
$ws.Cells.Item(47, 4).Value = 0
$ws.Cells.Item(47, 5).Value = 0
$ws.Cells.Item(47, 6).Value = 0
$ws.Cells.Item(48, 4).Value = 9757600
$ws.Cells.Item(48, 5).Value = 8244200
$ws.Cells.Item(48, 6).Value = 8208100
$ws.Cells.Item(49, 4).Value = 0
$ws.Cells.Item(49, 5).Value = 0
$ws.Cells.Item(49, 6).Value = 0
$ws.Cells.Item(50, 4).Value = 0
$ws.Cells.Item(50, 5).Value = 0
$ws.Cells.Item(50, 6).Value = 0
$ws.Cells.Item(51, 4).Value = 0
$ws.Cells.Item(51, 5).Value = 0
$ws.Cells.Item(51, 6).Value = 0
$ws.Cells.Item(52, 4).Value = 415200
$ws.Cells.Item(52, 5).Value = 401200
$ws.Cells.Item(52, 6).Value = 422000
$ws.Cells.Item(53, 4).Value = 0
$ws.Cells.Item(53, 5).Value = 0
$ws.Cells.Item(53, 6).Value = 0
$ws.Cells.Item(54, 4).Value = 11052600
$ws.Cells.Item(54, 5).Value = 10026600
$ws.Cells.Item(54, 6).Value = 9945400
$ws.Cells.Item(57, 4).Value = 592200
$ws.Cells.Item(57, 5).Value = 578000
$ws.Cells.Item(57, 6).Value = 563200
$ws.Cells.Item(58, 4).Value = 10600
$ws.Cells.Item(58, 5).Value = 10500
$ws.Cells.Item(58, 6).Value = 9700
$ws.Cells.Item(59, 4).Value = 243200
$ws.Cells.Item(59, 5).Value = 287100
$ws.Cells.Item(59, 6).Value = 330700
$ws.Cells.Item(60, 4).Value = 846100
$ws.Cells.Item(60, 5).Value = 875500
$ws.Cells.Item(60, 6).Value = 903600
$ws.Cells.Item(61, 4).Value = 3227100
$ws.Cells.Item(61, 5).Value = 2903900
$ws.Cells.Item(61, 6).Value = 2897300
$ws.Cells.Item(62, 4).Value = 1781800
$ws.Cells.Item(62, 5).Value = 1480300
$ws.Cells.Item(62, 6).Value = 1472900
$ws.Cells.Item(63, 4).Value = 0
$ws.Cells.Item(63, 5).Value = 0
$ws.Cells.Item(63, 6).Value = 0
$ws.Cells.Item(64, 4).Value = 0
$ws.Cells.Item(64, 5).Value = 0
$ws.Cells.Item(64, 6).Value = 0
$ws.Cells.Item(65, 4).Value = 0
$ws.Cells.Item(65, 5).Value = 0
$ws.Cells.Item(65, 6).Value = 0
$ws.Cells.Item(66, 4).Value = 6223300
$ws.Cells.Item(66, 5).Value = 5259700
$ws.Cells.Item(66, 6).Value = 5273800
$ws.Cells.Item(68, 4).Value = 0
$ws.Cells.Item(68, 5).Value = 0
$ws.Cells.Item(68, 6).Value = 0
$ws.Cells.Item(69, 4).Value = 0
$ws.Cells.Item(69, 5).Value = 0
$ws.Cells.Item(69, 6).Value = 0
$ws.Cells.Item(70, 4).Value = 0
$ws.Cells.Item(70, 5).Value = 0
$ws.Cells.Item(70, 6).Value = 0
$ws.Cells.Item(71, 4).Value = 0
$ws.Cells.Item(71, 5).Value = 0
$ws.Cells.Item(71, 6).Value = 0
$ws.Cells.Item(72, 4).Value = 5513500
$ws.Cells.Item(72, 5).Value = 5453400
$ws.Cells.Item(72, 6).Value = 5402700
$ws.Cells.Item(73, 4).Value = 0
$ws.Cells.Item(73, 5).Value = 0
$ws.Cells.Item(73, 6).Value = 0
$ws.Cells.Item(74, 4).Value = 0
$ws.Cells.Item(74, 5).Value = 0
$ws.Cells.Item(74, 6).Value = 0
$ws.Cells.Item(75, 4).Value = 0
$ws.Cells.Item(75, 5).Value = 0
$ws.Cells.Item(75, 6).Value = 0
$ws.Cells.Item(76, 4).Value = 4829300
$ws.Cells.Item(76, 5).Value = 4766900
$ws.Cells.Item(76, 6).Value = 4671600
$ws.Cells.Item(77, 4).Value = 0
$ws.Cells.Item(77, 5).Value = 0
$ws.Cells.Item(77, 6).Value = 0
$ws.Cells.Item(80, 4).Value = 43465
$ws.Cells.Item(80, 5).Value = 43373
$ws.Cells.Item(80, 6).Value = 43281
$ws.Cells.Item(81, 4).Value = 103400
$ws.Cells.Item(81, 5).Value = 93900
$ws.Cells.Item(81, 6).Value = 45500
$ws.Cells.Item(83, 4).Value = 261300
$ws.Cells.Item(83, 5).Value = 241800
$ws.Cells.Item(83, 6).Value = 238000
$ws.Cells.Item(84, 4).Value = 0
$ws.Cells.Item(84, 5).Value = 0
$ws.Cells.Item(84, 6).Value = 0
$ws.Cells.Item(85, 4).Value = 0
$ws.Cells.Item(85, 5).Value = 0
$ws.Cells.Item(85, 6).Value = 0
$ws.Cells.Item(86, 4).Value = 0
$ws.Cells.Item(86, 5).Value = 0
$ws.Cells.Item(86, 6).Value = 0
$ws.Cells.Item(87, 4).Value = 0
$ws.Cells.Item(87, 5).Value = 0
$ws.Cells.Item(87, 6).Value = 0
$ws.Cells.Item(88, 4).Value = 0
$ws.Cells.Item(88, 5).Value = 0
$ws.Cells.Item(88, 6).Value = 0
$ws.Cells.Item(89, 4).Value = 222500
$ws.Cells.Item(89, 5).Value = 372400
$ws.Cells.Item(89, 6).Value = 346000
$ws.Cells.Item(91, 4).Value = 63700
$ws.Cells.Item(91, 5).Value = -243200
$ws.Cells.Item(91, 6).Value = -341200
$ws.Cells.Item(92, 4).Value = 0
$ws.Cells.Item(92, 5).Value = 0
$ws.Cells.Item(92, 6).Value = 0
$ws.Cells.Item(93, 4).Value = 0
$ws.Cells.Item(93, 5).Value = 0
$ws.Cells.Item(93, 6).Value = 0
$ws.Cells.Item(94, 4).Value = -1038800
$ws.Cells.Item(94, 5).Value = -242700
$ws.Cells.Item(94, 6).Value = -340900
$ws.Cells.Item(96, 4).Value = -43300
$ws.Cells.Item(96, 5).Value = -43300
$ws.Cells.Item(96, 6).Value = -43300
$ws.Cells.Item(97, 4).Value = 0
$ws.Cells.Item(97, 5).Value = 0
$ws.Cells.Item(97, 6).Value = 0
$ws.Cells.Item(98, 4).Value = 0
$ws.Cells.Item(98, 5).Value = 0
$ws.Cells.Item(98, 6).Value = 0
$ws.Cells.Item(99, 4).Value = 0
$ws.Cells.Item(99, 5).Value = 0
$ws.Cells.Item(99, 6).Value = 0
$ws.Cells.Item(100, 4).Value = 271600
$ws.Cells.Item(100, 5).Value = -45800
$ws.Cells.Item(100, 6).Value = -45800
$ws.Cells.Item(101, 4).Value = -15600
$ws.Cells.Item(101, 5).Value = -37500
$ws.Cells.Item(101, 6).Value = 3300
$ws.Cells.Item(102, 4).Value = -560400
$ws.Cells.Item(102, 5).Value = 46400
$ws.Cells.Item(102, 6).Value = -37300

Write-Output "done"
